$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 17, shifting existing rows 17-39 down to 18-40
$ws.Rows.Item(17).Insert()

# Populate the new row 17 with data
$ws.Cells.Item(17, 1).Value = 6
$ws.Cells.Item(17, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(17, 3).Value = "Metropolitana"
$ws.Cells.Item(17, 4).Value = 44792
$ws.Cells.Item(17, 5).Value = 13
$ws.Cells.Item(17, 6).Value = 100112035
$ws.Cells.Item(17, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(17, 8).Value = "Sin especificar"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 250
$ws.Cells.Item(17, 11).Value = 16000
$ws.Cells.Item(17, 12).Value = 18000
$ws.Cells.Item(17, 13).Value = 16800
$ws.Cells.Item(17, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(17, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(17, 16).Value = 1120
$ws.Cells.Item(17, 17).Value = 15
$ws.Cells.Item(17, 18).Value = "Hortaliza"
